$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Sheet1 "Forecast Comparison": shift Week_Start_Date by one week and update MyForecast values ---
$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "2025-01-12"
$ws1.Range("D2").Value = 63
$ws1.Range("B3").NumberFormat = "@"
$ws1.Range("B3").Value = "2025-01-19"
$ws1.Range("D3").Value = 60
$ws1.Range("B4").NumberFormat = "@"
$ws1.Range("B4").Value = "2025-01-26"
$ws1.Range("D4").Value = 54
$ws1.Range("B5").NumberFormat = "@"
$ws1.Range("B5").Value = "2025-02-02"
$ws1.Range("D5").Value = 47
$ws1.Range("B6").NumberFormat = "@"
$ws1.Range("B6").Value = "2025-02-09"
$ws1.Range("D6").Value = 40
$ws1.Range("B7").NumberFormat = "@"
$ws1.Range("B7").Value = "2025-02-16"
$ws1.Range("D7").Value = 41
$ws1.Range("B8").NumberFormat = "@"
$ws1.Range("B8").Value = "2025-02-23"
$ws1.Range("D8").Value = 51
$ws1.Range("B9").NumberFormat = "@"
$ws1.Range("B9").Value = "2025-03-02"
$ws1.Range("D9").Value = 45
$ws1.Range("B10").NumberFormat = "@"
$ws1.Range("B10").Value = "2025-03-09"
$ws1.Range("D10").Value = 44
$ws1.Range("B11").NumberFormat = "@"
$ws1.Range("B11").Value = "2025-03-16"
$ws1.Range("D11").Value = 44
$ws1.Range("B12").NumberFormat = "@"
$ws1.Range("B12").Value = "2025-03-23"
$ws1.Range("D12").Value = 43
$ws1.Range("B13").NumberFormat = "@"
$ws1.Range("B13").Value = "2025-03-30"
$ws1.Range("D13").Value = 61
$ws1.Range("B14").NumberFormat = "@"
$ws1.Range("B14").Value = "2025-04-06"
$ws1.Range("D14").Value = 59
$ws1.Range("B15").NumberFormat = "@"
$ws1.Range("B15").Value = "2025-04-13"
$ws1.Range("D15").Value = 60
$ws1.Range("B16").NumberFormat = "@"
$ws1.Range("B16").Value = "2025-04-20"
$ws1.Range("D16").Value = 60
$ws1.Range("B17").NumberFormat = "@"
$ws1.Range("B17").Value = "2025-04-27"
$ws1.Range("D17").Value = 56

# --- Sheet2 "Summary": update summary metrics ---
$ws2.Range("B2").Value = "2022-12-25 to 2025-01-05"
$ws2.Range("B4").NumberFormat = "@"
$ws2.Range("B4").Value = "279"
$ws2.Range("B5").NumberFormat = "@"
$ws2.Range("B5").Value = "98"
$ws2.Range("B8").Value = "10684 units"
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "827"
$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "401"
$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "224"
$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "63"
$ws2.Range("B13").NumberFormat = "@"
$ws2.Range("B13").Value = "2025-01-12"
$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "40"
$ws2.Range("B15").NumberFormat = "@"
$ws2.Range("B15").Value = "2025-02-09"
